$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.072.24'
$ws.Cells.Item(2, 5).Value = '  +6.02%  '
$ws.Cells.Item(3, 4).Value = '3.110.05'
$ws.Cells.Item(3, 5).Value = '  +3.92%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '584.58'
$ws.Cells.Item(5, 5).Value = '  +3.71%  '
$ws.Cells.Item(6, 5).Value = '  +3.86%  '
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
$ws.Cells.Item(8, 4).Value = '3.102.68'
$ws.Cells.Item(8, 5).Value = '  +4.00%  '
$ws.Cells.Item(9, 5).Value = '  +1.89%  '
$ws.Cells.Item(10, 4).Value = '0.150'
$ws.Cells.Item(10, 5).Value = '  +13.07%  '
$ws.Cells.Item(11, 4).Value = '5.78'
$ws.Cells.Item(11, 5).Value = '  +8.73%  '
$ws.Cells.Item(12, 5).Value = '  +2.87%  '
$ws.Cells.Item(13, 5).Value = '  +7.47%  '
$ws.Cells.Item(14, 4).Value = '35.52'
$ws.Cells.Item(14, 5).Value = '  +5.14%  '
$ws.Cells.Item(15, 5).Value = '  +0.52%  '
$ws.Cells.Item(16, 4).Value = '3.623.85'
$ws.Cells.Item(16, 5).Value = '  +3.82%  '
$ws.Cells.Item(17, 5).Value = '  -0.28%  '
$ws.Cells.Item(18, 4).Value = '62.977.47'
$ws.Cells.Item(18, 5).Value = '  +5.88%  '
$ws.Cells.Item(19, 4).Value = '3.107.63'
$ws.Cells.Item(19, 5).Value = '  +3.88%  '
$ws.Cells.Item(20, 4).Value = '465.08'
$ws.Cells.Item(20, 5).Value = '  +6.97%  '
$ws.Cells.Item(21, 4).Value = '14.18'
$ws.Cells.Item(21, 5).Value = '  +4.33%  '
$ws.Cells.Item(22, 5).Value = '  +0.95%  '
$ws.Cells.Item(23, 5).Value = '  +6.89%  '
$ws.Cells.Item(24, 4).Value = '13.29'
$ws.Cells.Item(24, 5).Value = '  -0.90%  '
$ws.Cells.Item(25, 4).Value = '81.91'
$ws.Cells.Item(25, 5).Value = '  +2.34%  '
$ws.Cells.Item(26, 5).Value = '  -0.14%  '
$ws.Cells.Item(27, 4).Value = '8.38'
$ws.Cells.Item(27, 5).Value = '  +7.86%  '
$ws.Cells.Item(29, 5).Value = '  +4.99%  '
$ws.Cells.Item(31, 4).Value = '6.85'
$ws.Cells.Item(31, 5).Value = '  +9.63%  '
$ws.Cells.Item(32, 5).Value = '  +4.51%  '
$ws.Cells.Item(33, 5).Value = '  +2.47%  '
$ws.Cells.Item(34, 4).Value = '0.0₃0861'
$ws.Cells.Item(34, 5).Value = '  +10.47%  '
$ws.Cells.Item(35, 5).Value = '  +15.67%  '
$ws.Cells.Item(36, 5).Value = '  +4.42%  '
$ws.Cells.Item(37, 4).Value = '3.32'
$ws.Cells.Item(37, 5).Value = '  +19.24%  '
$ws.Cells.Item(38, 4).Value = '6.03'
$ws.Cells.Item(38, 5).Value = '  +2.51%  '
$ws.Cells.Item(39, 4).Value = '50.82'
$ws.Cells.Item(39, 5).Value = '  +4.06%  '
$ws.Cells.Item(40, 4).Value = '433.07'
$ws.Cells.Item(40, 5).Value = '  +7.82%  '
$ws.Cells.Item(41, 5).Value = '  +0.92%  '
$ws.Cells.Item(42, 4).Value = '2.922.77'
$ws.Cells.Item(42, 5).Value = '  +5.93%  '
$ws.Cells.Item(43, 4).Value = '0.0369'
$ws.Cells.Item(43, 5).Value = '  +4.21%  '
$ws.Cells.Item(44, 5).Value = '  +11.63%  '
$ws.Cells.Item(45, 5).Value = '  +5.46%  '
$ws.Cells.Item(46, 4).Value = '2.16'
$ws.Cells.Item(46, 5).Value = '  +7.44%  '
$ws.Cells.Item(47, 4).Value = '35.35'
$ws.Cells.Item(47, 5).Value = '  +1.76%  '
$ws.Cells.Item(49, 5).Value = '  -0.07%  '
$ws.Cells.Item(50, 5).Value = '  +0.70%  '
$ws.Cells.Item(51, 4).Value = '24.48'
$ws.Cells.Item(51, 5).Value = '  +4.18%  '
